# Apply updated SDV generative model evaluation metrics to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new value
$updates = @{
    "G5"  = 0.95
    "I5"  = 0.9737174218898829
    "G6"  = 0.9264285714285714
    "I6"  = 0.9570521398934787
    "G7"  = 0.9478571428571428
    "I7"  = 0.9839147155248877
    "G8"  = 0.9764285714285714
    "I8"  = 0.9595425914428009
    "G9"  = 0.9314285714285714
    "I9"  = 0.9757416146699933
    "G11" = 0.9507142857142857
    "I11" = 0.9692530767779589
    "G12" = 0.9314285714285714
    "I12" = 0.9722338959612816
    "G13" = 0.9557142857142857
    "I13" = 0.9815145730838175
    "G16" = 0.9614285714285714
    "I16" = 0.9831423389731938
    "G17" = 0.9607142857142857
    "I17" = 0.9974540276653338
    "G18" = 0.9392857142857143
    "I18" = 0.9805388373224927
    "E19" = 0.921994884910486
    "G19" = 0.9357142857142857
    "I19" = 0.8933328783999316
    "E20" = 0.984251968503937
    "G20" = 0.9485714285714286
    "I20" = 0.9432262158276035
    "E21" = 0.8539325842696632
    "G21" = 0.9257142857142857
    "I21" = 0.9300671349943236
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
